# Append the "final round" of USDA families resazurin metadata rows
# (rows 152-181) for sampling date 2025-06-10, following the exact same
# pattern used for the previous sampling dates already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = 20250610

# Each entry: row, sample, well, type, family (family omitted for blanks)
$rows = @(
    @(152, "A1", "A01", "sample", "A"),
    @(153, "A2", "A02", "sample", "A"),
    @(154, "A3", "A03", "sample", "A"),
    @(155, "A4", "A04", "sample", "A"),
    @(156, "A5", "A05", "sample", "A"),
    @(157, "B1", "B01", "sample", "B"),
    @(158, "B2", "B02", "sample", "B"),
    @(159, "B3", "B03", "sample", "B"),
    @(160, "B4", "B04", "sample", "B"),
    @(161, "B5", "B05", "sample", "B"),
    @(162, "C1", "C01", "sample", "C"),
    @(163, "C2", "C02", "sample", "C"),
    @(164, "C3", "C03", "sample", "C"),
    @(165, "C4", "C04", "sample", "C"),
    @(166, "C5", "C05", "sample", "C"),
    @(167, "D1", "D01", "sample", "D"),
    @(168, "D2", "D02", "sample", "D"),
    @(169, "D3", "D03", "sample", "D"),
    @(170, "D4", "D04", "sample", "D"),
    @(171, "D5", "D05", "sample", "D"),
    @(172, "E1", "E01", "sample", "E"),
    @(173, "E2", "E02", "sample", "E"),
    @(174, "E3", "E03", "sample", "E"),
    @(175, "E4", "E04", "sample", "E"),
    @(176, "E5", "E05", "sample", "E"),
    @(177, "Blank1", "H01", "blank", $null),
    @(178, "Blank2", "H02", "blank", $null),
    @(179, "Blank3", "H03", "blank", $null),
    @(180, "Blank4", "H04", "blank", $null),
    @(181, "Blank5", "H05", "blank", $null)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $sample = $r[1]
    $well = $r[2]
    $type = $r[3]
    $family = $r[4]

    $ws.Cells.Item($rowNum, 1).Value = $date
    $ws.Cells.Item($rowNum, 2).Value = $sample
    $ws.Cells.Item($rowNum, 3).Value = $well
    $ws.Cells.Item($rowNum, 4).Value = $type
    if ($family) {
        $ws.Cells.Item($rowNum, 5).Value = $family
    }
}

# Match the existing bold formatting used on the "sample" and "well" columns
# (B, C) for every other data row by copying the style from an existing row.
$fmtSource = $ws.Range("B2:C2")
$fmtSource.Copy()
$fmtTarget = $ws.Range("B152:C181")
$fmtTarget.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Restore the selection/scroll position left behind by the edit session
$ws.Range("B178").Select()
